# "Fix template tests - fix input data and ref value"
#
# The "ITR input data" sheet had a 2019_production figure (AV2) that was
# off by a factor of 10 (a data-entry error) - correct it.
#
# Reflect the editor's resulting UI state: the "ITR input data" sheet
# becomes the active tab (was "Portfolio"), with the just-fixed cell AV2
# selected; the "Portfolio" sheet's previous special view state is
# cleared since it's no longer the active sheet.

$wb = $excel.ActiveWorkbook

$inputSheet = $wb.Worksheets.Item("ITR input data")

# Fix the mis-entered 2019 production value (was 10x too large).
$inputSheet.Range("AV2").Value = 75904.354999999996

# Make "ITR input data" the active sheet with AV2 selected, matching the
# saved workbook view (activeTab + per-sheet tabSelected/selection).
[void]$inputSheet.Activate()
[void]$inputSheet.Range("AV2").Select()
